$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Object Type 2")

# Row 21 ("New Technologies"): parameter_index changes from InvCountry -> Year
$ws2.Range("C21").Value = "Year"

# Row 27 ("VRE Capacities"): parameter_index changes from Bus -> Year
$ws2.Range("C27").Value = "Year"

# Make "Object Type 2" the active tab / sheet, with selection on C28
$ws2.Activate() | Out-Null
$ws2.Range("C28").Select() | Out-Null
